# Updated cryptos list on Sun May 12 09:48:32 UTC 2024 with GitHub Actions
# Refresh price / 1h-volume figures and fix the Hedera / EthereumClassic
# ranking swap (rows 31-32). Values are written as literal text (matching
# the sheet's existing inlineStr cells) without disturbing cell styling:
# NumberFormat is forced to Text ("@") just long enough for the assignment
# to stick as a string (Excel would otherwise re-parse numeric-looking
# strings like "144.90" into the number 144.9), then the cell style is put
# back to "Normal" so no formatting/style index drifts from the original.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue $ws "D2" "60.996.42"
Set-TextValue $ws "E2" "  +0.34%  "
Set-TextValue $ws "D3" "2.912.70"
Set-TextValue $ws "E3" "  +0.22%  "
Set-TextValue $ws "E4" "  +0.02%  "
Set-TextValue $ws "D5" "590.02"
Set-TextValue $ws "E5" "  +0.80%  "
Set-TextValue $ws "D6" "144.90"
Set-TextValue $ws "E6" "  +0.11%  "
Set-TextValue $ws "E8" "  +0.25%  "
Set-TextValue $ws "E9" "  +4.01%  "
Set-TextValue $ws "E10" "  -2.21%  "
Set-TextValue $ws "E11" "  -1.58%  "
Set-TextValue $ws "E12" "  -0.66%  "
Set-TextValue $ws "E13" "  +0.06%  "
Set-TextValue $ws "E14" "  -0.25%  "
Set-TextValue $ws "D15" "3.395.96"
Set-TextValue $ws "E15" "  +0.26%  "
Set-TextValue $ws "D16" "60.912.28"
Set-TextValue $ws "E16" "  +0.29%  "
Set-TextValue $ws "E17" "  -0.30%  "
Set-TextValue $ws "D18" "2.914.41"
Set-TextValue $ws "D19" "433.18"
Set-TextValue $ws "E19" "  +1.15%  "
Set-TextValue $ws "E20" "  -1.54%  "
Set-TextValue $ws "E21" "  -1.04%  "
Set-TextValue $ws "E22" "  -0.22%  "
Set-TextValue $ws "D23" "81.42"
Set-TextValue $ws "E23" "  +0.97%  "
Set-TextValue $ws "E24" "  -0.46%  "
Set-TextValue $ws "E25" "  -1.72%  "
Set-TextValue $ws "E26" "  -1.18%  "
Set-TextValue $ws "E27" "  +0.02%  "
Set-TextValue $ws "D28" "2.28"
Set-TextValue $ws "E28" "  +4.47%  "
Set-TextValue $ws "E29" "  -0.66%  "
Set-TextValue $ws "E30" "  -4.31%  "
Set-TextValue $ws "B31" "Hedera"
Set-TextValue $ws "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D31" "0.109"
Set-TextValue $ws "E31" "  +2.47%  "
Set-TextValue $ws "B32" "EthereumClassic"
Set-TextValue $ws "C32" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D32" "26.44"
Set-TextValue $ws "E32" "  -0.18%  "
Set-TextValue $ws "E33" "  +0.08%  "
Set-TextValue $ws "E34" "  -0.86%  "
Set-TextValue $ws "E35" "  +0.12%  "
Set-TextValue $ws "D36" "5.60"
Set-TextValue $ws "E36" "  -0.01%  "
Set-TextValue $ws "E37" "  -1.06%  "
Set-TextValue $ws "E38" "  -1.09%  "
Set-TextValue $ws "E39" "  -2.67%  "
Set-TextValue $ws "D40" "8.56"
Set-TextValue $ws "E40" "  -0.55%  "
Set-TextValue $ws "D41" "41.73"
Set-TextValue $ws "E41" "  +0.66%  "
Set-TextValue $ws "D42" "0.286"
Set-TextValue $ws "E42" "  -3.28%  "
Set-TextValue $ws "D43" "376.09"
Set-TextValue $ws "E43" "  -0.51%  "
Set-TextValue $ws "E44" "  -1.44%  "
Set-TextValue $ws "D45" "2.685.38"
Set-TextValue $ws "E45" "  -0.55%  "
Set-TextValue $ws "D46" "132.74"
Set-TextValue $ws "E46" "  +0.29%  "
Set-TextValue $ws "D48" "23.79"
Set-TextValue $ws "E48" "  -1.84%  "
Set-TextValue $ws "E50" "  -2.92%  "
Set-TextValue $ws "D51" "0.123"
Set-TextValue $ws "E51" "  -0.97%  "
